{"js": "// Replace the date line and the twenty-five \"NNN\u00d7N=\" multiplication\n// prompts with the next day's values (diff: 2026-02-22 Sunday ->\n// 2026-02-23 Monday, plus per-cell operand/operator updates).\nconst replacements = [\n  [\"2026-02-22 Sunday\", \"2026-02-23 Monday\"],\n  [\"675\u00d77=\", \"636\u00d73=\"],\n  [\"219\u00d74=\", \"634\u00d76=\"],\n  [\"778\u00d79=\", \"146\u00d73=\"],\n  [\"778\u00d74=\", \"313\u00d77=\"],\n  [\"749\u00d76=\", \"497\u00d79=\"],\n  [\"848\u00d74=\", \"243\u00d77=\"],\n  [\"429\u00d77=\", \"661\u00d79=\"],\n  [\"647\u00d72=\", \"386\u00d76=\"],\n  [\"189\u00d72=\", \"716\u00d75=\"],\n  [\"155\u00d76=\", \"326\u00d77=\"],\n  [\"737\u00d76=\", \"645\u00d74=\"],\n  [\"970\u00d73=\", \"867\u00d74=\"],\n  [\"981\u00d79=\", \"525\u00d79=\"],\n  [\"942\u00d76=\", \"794\u00d77=\"],\n  [\"965\u00d75=\", \"682\u00d77=\"],\n  [\"843\u00d77=\", \"911\u00d76=\"],\n  [\"402\u00d79=\", \"619\u00d76=\"],\n  [\"417\u00d77=\", \"671\u00d74=\"],\n  [\"242\u00d72=\", \"356\u00d79=\"],\n  [\"169\u00d73=\", \"890\u00d76=\"],\n  [\"547\u00d76=\", \"703\u00d77=\"],\n  [\"949\u00d73=\", \"346\u00d77=\"],\n  [\"137\u00d72=\", \"946\u00d77=\"],\n  [\"628\u00d73=\", \"125\u00d78=\"],\n  [\"309\u00d77=\", \"393\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the twenty-five \"NNN\u00d7N=\" multiplication\n# prompts with the next day's values (diff: 2026-02-22 Sunday ->\n# 2026-02-23 Monday, plus per-cell operand/operator updates).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-22 Sunday\", \"2026-02-23 Monday\"),\n    @(\"675\u00d77=\", \"636\u00d73=\"),\n    @(\"219\u00d74=\", \"634\u00d76=\"),\n    @(\"778\u00d79=\", \"146\u00d73=\"),\n    @(\"778\u00d74=\", \"313\u00d77=\"),\n    @(\"749\u00d76=\", \"497\u00d79=\"),\n    @(\"848\u00d74=\", \"243\u00d77=\"),\n    @(\"429\u00d77=\", \"661\u00d79=\"),\n    @(\"647\u00d72=\", \"386\u00d76=\"),\n    @(\"189\u00d72=\", \"716\u00d75=\"),\n    @(\"155\u00d76=\", \"326\u00d77=\"),\n    @(\"737\u00d76=\", \"645\u00d74=\"),\n    @(\"970\u00d73=\", \"867\u00d74=\"),\n    @(\"981\u00d79=\", \"525\u00d79=\"),\n    @(\"942\u00d76=\", \"794\u00d77=\"),\n    @(\"965\u00d75=\", \"682\u00d77=\"),\n    @(\"843\u00d77=\", \"911\u00d76=\"),\n    @(\"402\u00d79=\", \"619\u00d76=\"),\n    @(\"417\u00d77=\", \"671\u00d74=\"),\n    @(\"242\u00d72=\", \"356\u00d79=\"),\n    @(\"169\u00d73=\", \"890\u00d76=\"),\n    @(\"547\u00d76=\", \"703\u00d77=\"),\n    @(\"949\u00d73=\", \"346\u00d77=\"),\n    @(\"137\u00d72=\", \"946\u00d77=\"),\n    @(\"628\u00d73=\", \"125\u00d78=\"),\n    @(\"309\u00d77=\", \"393\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
